# Addition of staging work
#
# A batch of new "product" entries (prodbJjt, prodUebe, prodtVOI, prodTSDr,
# prodLCir, prodmNla, prodiIIH, prodhsKd, prodhvAf, prodtmeE, prodUruY,
# prodHtae, prodptjn, produofr, prodfZRT, prodCgnG, prodBtxr, prodmfTT,
# prodFTCZ, prodtXgH, prodJIOl) was staged into the workbook. The six
# visible "Sku" cells on the Input sheet (B2:B7) are rotated to point at
# freshly staged product codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "prodJIOl"
$ws.Range("B3").Value = "prodbJjt"
$ws.Range("B4").Value = "prodmfTT"
$ws.Range("B5").Value = "prodUebe"
$ws.Range("B6").Value = "prodFTCZ"
$ws.Range("B7").Value = "prodtXgH"
